$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E9").Value = 13.021
$ws.Range("E13").Value = 12.734
$ws.Range("E16").Value = 12.963
$ws.Range("E18").Value = 13.115
$ws.Range("E20").Value = 13.127
